$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.418.11"
$ws.Range("E2").Value = "  +5.24%  "

$ws.Range("D3").Value = "2.455.89"
$ws.Range("E3").Value = "  +3.51%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.518"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.94%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"

$ws.Range("E9").Value = "  +10.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.65%  "

$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("E12").Value = "  -2.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.39%  "

$ws.Range("D15").Value = "2.835.25"
$ws.Range("E15").Value = "  +3.60%  "

$ws.Range("D16").Value = "2.444.92"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("E17").Value = "  +4.71%  "

$ws.Range("D18").Value = "45.294.20"
$ws.Range("E18").Value = "  +4.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.42%  "

$ws.Range("D21").Value = "0.0₃0931"
$ws.Range("E21").Value = "  +4.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.87%  "

$ws.Range("E24").Value = "  +3.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.21%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  +4.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.62%  "

$ws.Range("E32").Value = "  +15.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.62%  "

$ws.Range("E34").Value = "  +3.96%  "

$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0766"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.23%  "

$ws.Range("E41").Value = "  +2.43%  "

$ws.Range("E42").Value = "  -2.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0292"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.63%  "

$ws.Range("D45").Value = "1.953.09"
$ws.Range("E45").Value = "  +1.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.91%  "

$ws.Range("E47").Value = "  -0.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("E49").Value = "  +17.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.48%  "
